$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.045.89"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.92%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.737.00"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.58%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "601.31"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.98%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.40"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.06%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.738.33"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.77%  "
$ws.Range("E9").Value = "  +0.51%  "
$ws.Range("E10").Value = "  +1.54%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.41"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.94%  "
$ws.Range("E12").Value = "  +0.10%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "37.83"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.55%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000246"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.90%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.360.42"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.51%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.733.01"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.34%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "69.079.57"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.08%  "
$ws.Range("E18").Value = "  +0.41%  "
$ws.Range("E19").Value = "  -1.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.05"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.33%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.80"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +17.22%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "492.17"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.73%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.723"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.51%  "
$ws.Range("E24").Value = "  +6.38%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.59"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.08%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.30"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.26"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.64%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.11"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.04%  "
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("E30").Value = "  +2.60%  "
$ws.Range("E31").Value = "  +6.72%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.01"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.03%  "
$ws.Range("E33").Value = "  +0.93%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.881.22"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.54%  "
$ws.Range("E35").Value = "  +0.49%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.671.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.28%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.01"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.23%  "
$ws.Range("E39").Value = "  +0.47%  "
$ws.Range("E40").Value = "  +2.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.323"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.05%  "
$ws.Range("B42").Value = "Bittensor"
$ws.Range("C42").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "432.02"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.01%  "
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.95"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.68%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "48.55"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.12%  "
$ws.Range("E45").Value = "  +2.13%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.45"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.88%  "
$ws.Range("E47").Value = "  +0.02%  "
$ws.Range("E48").Value = "  -0.04%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "141.72"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.39%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.768.01"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.20%  "
$ws.Range("E51").Value = "  +1.00%  "
